$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 98.912777
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.32975033333333
$ws.Range("N2").Value = 48.989251
$ws.Range("O2").Value = 0.3418592381614251
$ws.Range("P2").Value = 0.3418592381614251
$ws.Range("Q2").Value = 1615.220953186676
$ws.Range("R2").Value = 14536.98857868008
$ws.Range("S2").Value = 0.277617909265453
$ws.Range("T2").Value = 0.277617909265453

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 98.912777
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 26.81766766666667
$ws.Range("N3").Value = 80.453003
$ws.Range("O3").Value = 0.5614211638667195
$ws.Range("P3").Value = 0.5614211638667195
$ws.Range("Q3").Value = 2652.60998157311
$ws.Range("R3").Value = 23873.48983415799
$ws.Range("S3").Value = 0.4559203096815507
$ws.Range("T3").Value = 0.4559203096815507

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 98.912777
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.620050333333333
$ws.Range("N4").Value = 13.860151
$ws.Range("O4").Value = 0.09671959797185539
$ws.Range("P4").Value = 0.09671959797185539
$ws.Range("Q4").Value = 456.9820083497757
$ws.Range("R4").Value = 4112.838075147981
$ws.Range("S4").Value = 0.07854429419064761
$ws.Range("T4").Value = 0.07854429419064761

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.32975033333333
$ws.Range("N5").Value = 48.989251
$ws.Range("O5").Value = 0.3418592381614251
$ws.Range("P5").Value = 0.3418592381614251
$ws.Range("Q5").Value = 278.3997516738743
$ws.Range("R5").Value = 2505.597765064868
$ws.Range("S5").Value = 0.04785026893518127
$ws.Range("T5").Value = 0.04785026893518127

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 26.81766766666667
$ws.Range("N6").Value = 80.453003
$ws.Range("O6").Value = 0.5614211638667195
$ws.Range("P6").Value = 0.5614211638667195
$ws.Range("Q6").Value = 457.2042968490672
$ws.Range("R6").Value = 4114.838671641604
$ws.Range("S6").Value = 0.07858250027527355
$ws.Range("T6").Value = 0.07858250027527355

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.620050333333333
$ws.Range("N7").Value = 13.860151
$ws.Range("O7").Value = 0.09671959797185539
$ws.Range("P7").Value = 0.09671959797185539
$ws.Range("Q7").Value = 78.76549483400757
$ws.Range("R7").Value = 708.8894535060681
$ws.Range("S7").Value = 0.01353790758777311
$ws.Range("T7").Value = 0.01353790758777311

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.32975033333333
$ws.Range("N8").Value = 48.989251
$ws.Range("O8").Value = 0.3418592381614251
$ws.Range("P8").Value = 0.3418592381614251
$ws.Range("Q8").Value = 95.36554598966268
$ws.Range("R8").Value = 858.289913906964
$ws.Range("S8").Value = 0.01639105996079085
$ws.Range("T8").Value = 0.01639105996079085

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 26.81766766666667
$ws.Range("N9").Value = 80.453003
$ws.Range("O9").Value = 0.5614211638667195
$ws.Range("P9").Value = 0.5614211638667195
$ws.Range("Q9").Value = 156.6148573613214
$ws.Range("R9").Value = 1409.533716251892
$ws.Range("S9").Value = 0.02691835390989518
$ws.Range("T9").Value = 0.02691835390989518

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.04794681006412069
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.620050333333333
$ws.Range("N10").Value = 13.860151
$ws.Range("O10").Value = 0.09671959797185539
$ws.Range("P10").Value = 0.09671959797185539
$ws.Range("Q10").Value = 26.98103850606267
$ws.Range("R10").Value = 242.829346554564
$ws.Range("S10").Value = 0.004637396193434664
$ws.Range("T10").Value = 0.004637396193434663
